$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to Text format before assigning, so that
# values like "1.00", "0.999", "224.97" are preserved exactly as text
# instead of being auto-coerced into numbers (which would drop trailing
# zeros / introduce floating point noise).
$dCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "34.130.71"
$ws.Range("E2").Value2 = "  -0.26%  "

$ws.Range("D3").Value2 = "1.825.30"
$ws.Range("E3").Value2 = "  +2.33%  "

$ws.Range("D4").Value2 = "1.00"
$ws.Range("E4").Value2 = "  +0.22%  "

$ws.Range("D5").Value2 = "224.97"
$ws.Range("E5").Value2 = "  +0.39%  "

$ws.Range("D6").Value2 = "0.558"
$ws.Range("E6").Value2 = "  +1.15%  "

$ws.Range("E7").Value2 = "  -0.03%  "

$ws.Range("D8").Value2 = "31.94"
$ws.Range("E8").Value2 = "  -0.19%  "

$ws.Range("D9").Value2 = "0.290"
$ws.Range("E9").Value2 = "  +3.58%  "

$ws.Range("D10").Value2 = "0.0722"
$ws.Range("E10").Value2 = "  +9.81%  "

$ws.Range("E11").Value2 = "  -0.06%  "

$ws.Range("D12").Value2 = "2.088.32"
$ws.Range("E12").Value2 = "  +2.38%  "

$ws.Range("D13").Value2 = "1.824.96"
$ws.Range("E13").Value2 = "  +2.51%  "

$ws.Range("D14").Value2 = "10.82"
$ws.Range("E14").Value2 = "  -3.20%  "

$ws.Range("D15").Value2 = "0.644"
$ws.Range("E15").Value2 = "  +2.92%  "

$ws.Range("D16").Value2 = "34.191.03"
$ws.Range("E16").Value2 = "  -0.11%  "

$ws.Range("E17").Value2 = "  +3.08%  "

$ws.Range("D18").Value2 = "69.68"
$ws.Range("E18").Value2 = "  +1.36%  "

$ws.Range("D19").Value2 = "250.83"
$ws.Range("E19").Value2 = "  -1.25%  "

$ws.Range("D20").Value2 = "0.0₃0791"
$ws.Range("E20").Value2 = "  +6.85%  "

$ws.Range("D21").Value2 = "11.15"
$ws.Range("E21").Value2 = "  +7.74%  "

$ws.Range("D22").Value2 = "0.999"
$ws.Range("E22").Value2 = "  +0.04%  "

$ws.Range("E23").Value2 = "  +1.65%  "

$ws.Range("E24").Value2 = "  +1.07%  "

$ws.Range("D25").Value2 = "160.62"
$ws.Range("E25").Value2 = "  +2.14%  "

$ws.Range("D26").Value2 = "16.62"
$ws.Range("E26").Value2 = "  +1.51%  "

$ws.Range("D27").Value2 = "7.25"
$ws.Range("E27").Value2 = "  +3.43%  "

$ws.Range("E28").Value2 = "  +0.93%  "

$ws.Range("D29").Value2 = "1.00"
$ws.Range("E29").Value2 = "  +0.05%  "

$ws.Range("D30").Value2 = "0.0533"
$ws.Range("E30").Value2 = "  +3.75%  "

$ws.Range("E32").Value2 = "  +2.10%  "

$ws.Range("D33").Value2 = "3.58"
$ws.Range("E33").Value2 = "  -0.29%  "

$ws.Range("D34").Value2 = "1.89"
$ws.Range("E34").Value2 = "  +1.39%  "

$ws.Range("D35").Value2 = "1.437.85"
$ws.Range("E35").Value2 = "  -0.28%  "

$ws.Range("B36").Value2 = "TrustWalletToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value2 = "1.07"
$ws.Range("E36").Value2 = "  +1.01%  "

$ws.Range("B37").Value2 = "ImmutableX"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value2 = "0.643"
$ws.Range("E37").Value2 = "  +3.12%  "

$ws.Range("E38").Value2 = "  +1.41%  "

$ws.Range("D39").Value2 = "0.960"
$ws.Range("E39").Value2 = "  +8.11%  "

$ws.Range("D40").Value2 = "81.64"
$ws.Range("E40").Value2 = "  -1.48%  "

$ws.Range("E41").Value2 = "  -3.23%  "

$ws.Range("E42").Value2 = "  +0.06%  "

$ws.Range("D43").Value2 = "2.15"
$ws.Range("E43").Value2 = "  +4.61%  "

$ws.Range("D44").Value2 = "6.09"
$ws.Range("E44").Value2 = "  +4.33%  "

$ws.Range("B45").Value2 = "WEMIXToken"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value2 = "1.06"
$ws.Range("E45").Value2 = "  +0.83%  "

$ws.Range("B46").Value2 = "RocketPoolETH"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value2 = "1.984.25"
$ws.Range("E46").Value2 = "  +2.13%  "

$ws.Range("D47").Value2 = "0.0497"
$ws.Range("E47").Value2 = "  -2.40%  "

$ws.Range("D48").Value2 = "106.98"
$ws.Range("E48").Value2 = "  +8.76%  "

$ws.Range("D49").Value2 = "0.999"
$ws.Range("E49").Value2 = "  +0.01%  "

$ws.Range("E50").Value2 = "  -2.54%  "

$ws.Range("E51").Value2 = "  +4.48%  "

